$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last row (69) previously carried the "latest entry" date-only format.
# Today's daily update appends row 70 as the new latest entry, so row 69
# reverts to the standard timestamp format used by all the other data rows.
$ws.Range("A69").NumberFormat = $ws.Range("A68").NumberFormat

# Append the new day's data as row 70.
$ws.Range("A70").Value = 45657
$ws.Range("B70").Value = 164
$ws.Range("C70").Value = 158
$ws.Range("D70").Value = 162

# New last row gets the date-only "latest" formatting that row 69 used to have.
$ws.Range("A70").NumberFormat = "YYYY-MM-DD"
